$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 475.6251530700926
$ws.Range("R2").Value = 4280.626377630833
$ws.Range("S2").Value = 0.001120500119079388
$ws.Range("T2").Value = 0.001120500119079388
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 2497.528577012537
$ws.Range("R3").Value = 22477.75719311283
$ws.Range("S3").Value = 0.005883795358346645
$ws.Range("T3").Value = 0.005883795358346645
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 1245.313647320055
$ws.Range("R4").Value = 11207.82282588049
$ws.Range("S4").Value = 0.002933768496275623
$ws.Range("T4").Value = 0.002933768496275623
$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 2956.923885395312
$ws.Range("R5").Value = 26612.31496855781
$ws.Range("S5").Value = 0.006966060445516146
$ws.Range("T5").Value = 0.006966060445516146
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 26648.52318584242
$ws.Range("R6").Value = 239836.7086725818
$ws.Range("S6").Value = 0.06277984503192553
$ws.Range("T6").Value = 0.06277984503192552
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.3296597247129734
$ws.Range("T7").Value = 0.3296597247129733
$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 69772.94911769287
$ws.Range("R8").Value = 627956.542059236
$ws.Range("S8").Value = 0.1643743971281804
$ws.Range("T8").Value = 0.1643743971281804
$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 165671.756865886
$ws.Range("R9").Value = 1491045.811792974
$ws.Range("S9").Value = 0.390297321531605
$ws.Range("T9").Value = 0.390297321531605
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 608.5742239108156
$ws.Range("R10").Value = 5477.16801519734
$ws.Range("S10").Value = 0.001433707796904978
$ws.Range("T10").Value = 0.001433707796904978
$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 3195.649989575923
$ws.Range("R11").Value = 28760.84990618331
$ws.Range("S11").Value = 0.007528462636474944
$ws.Range("T11").Value = 0.007528462636474944
$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 1593.409813487435
$ws.Range("R12").Value = 14340.68832138692
$ws.Range("S12").Value = 0.003753829826346089
$ws.Range("T12").Value = 0.003753829826346089
$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 3783.457723171787
$ws.Range("R13").Value = 34051.11950854608
$ws.Range("S13").Value = 0.008913247758200597
$ws.Range("T13").Value = 0.008913247758200597
$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 403.910916362877
$ws.Range("R14").Value = 3635.198247265893
$ws.Range("S14").Value = 0.0009515523452885424
$ws.Range("T14").Value = 0.0009515523452885423
$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 2120.953969049116
$ws.Range("R15").Value = 19088.58572144204
$ws.Range("S15").Value = 0.004996643174864235
$ws.Range("T15").Value = 0.004996643174864234
$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 1057.546627215726
$ws.Range("R16").Value = 9517.919644941534
$ws.Range("S16").Value = 0.002491418113777924
$ws.Range("T16").Value = 0.002491418113777924
$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 2511.082158830424
$ws.Range("R17").Value = 22599.73942947382
$ws.Range("S17").Value = 0.00591572552424066
$ws.Range("T17").Value = 0.005915725524240659
